# Adds three new template_type rows (for address line1/line2/line3) right
# after the existing "mosip.address.template.property" (Address) row, i.e.
# before the old row 1714 (Province). All rows from the old 1714 onward are
# pushed down by 3 rows, but keep their original content unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 1714 (this shifts existing rows 1714+
# down to 1717+, exactly like the rows below "Address" sliding down).
$null = $ws.Range("A1714:A1716").EntireRow.Insert()

# New row 1714: mosip.address.line1.template.property / Address line1
$ws.Cells.Item(1714, 1).Value = "eng"
$ws.Cells.Item(1714, 2).Value = "mosip.address.line1.template.property"
$ws.Cells.Item(1714, 3).Value = "Address line1"

# New row 1715: mosip.address.line2.template.property / Address line2
$ws.Cells.Item(1715, 1).Value = "eng"
$ws.Cells.Item(1715, 2).Value = "mosip.address.line2.template.property"
$ws.Cells.Item(1715, 3).Value = "Address line2"

# New row 1716: mosip.address.line3.template.property / Address line3
$ws.Cells.Item(1716, 1).Value = "eng"
$ws.Cells.Item(1716, 2).Value = "mosip.address.line3.template.property"
$ws.Cells.Item(1716, 3).Value = "Address line3"

# Column D holds the literal text "TRUE" (stored as a shared string, not a
# boolean) for every data row. Copy it down from the row right above (which
# already has the exact same style/shared-string) so the new cells match the
# existing "TRUE" text cells exactly instead of becoming real booleans.
$trueCell = $ws.Cells.Item(1713, 4)
for ($r = 1714; $r -le 1716; $r++) {
    $null = $trueCell.Copy()
    $null = $ws.Cells.Item($r, 4).PasteSpecial()
}
$excel.CutCopyMode = $false

# Restore the selection to mirror where the author ended up after inserting
# the new rows.
$null = $ws.Range("C1716").Select()
